$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Region: S-West Region -> North Region
$ws.Range("A5").Value = "North Region`t"

# Test date: placeholder text "yyyy-mm-dd" -> real date value (2017-01-25)
$ws.Range("B5").Value = 42760

# Test Centre: Shijiazhuang -> BJ-UIBE
$ws.Range("C5").Value = "BJ-UIBE"

# Candidate Number: 123456 (number) -> "005774" (text, so format must be applied before the value)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "005774"

# Candidate Name: "Hello Kitty" -> 990801 (numeric), displayed via text format
$ws.Range("E5").Value = 990801
$ws.Range("E5").NumberFormat = "@"

# Examiner Name, Listening/Reading original scores stay the same values (Bill Gates / 6.0 / 6.0)

# Move the active selection to C5 (matches the saved selection state in the workbook)
$ws.Range("C5").Select()
